# UseCaseCourseData.xlsx — "Fixed UseCase Data and Creation" edit
#
# Content change: the "Survey" sheet's AssessmentType column (BC) held a
# stray numeric 1 in row 2 and was otherwise empty for the data rows — it
# should instead carry the literal "SURVEY" marker (same shared string
# already used in column C) for every data row (rows 2-14), styled with a
# small accent font.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Survey")

# ---------------------------------------------------------------------
# 1. Data fix: BC2:BC14 = "SURVEY" (BC2 previously held a stray numeric 1,
#    BC3:BC14 previously had no AssessmentType value at all).
# ---------------------------------------------------------------------
$dataRange = $ws.Range("BC2:BC14")
$dataRange.Value = "SURVEY"

# ---------------------------------------------------------------------
# 2. New font/style for that column: 8pt accent-gold text, vertically
#    centred. Built once on a scratch cell (so only a single new font +
#    single new cell style get minted) and then format-painted onto the
#    whole BC2:BC14 block in one shot.
# ---------------------------------------------------------------------
$scratch = $ws.Range("ZZ1")
$scratch.Font.Size = 8
$scratch.Font.Color = 7185097          # RGB(201,162,109) == FFC9A26D
$scratch.Font.Name = "MesloLGM NF"
$scratch.Font.Family = 3
$scratch.VerticalAlignment = -4108     # xlVAlignCenter

$scratch.Copy()
$dataRange.PasteSpecial(-4122)         # xlPasteFormats
$excel.CutCopyMode = 0
$scratch.Clear()

# ---------------------------------------------------------------------
# 3. Leave the cursor where the author left it: Survey sheet, BD1
#    selected, scrolled over to the right-hand columns.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 26
$ws.Range("BD1").Select()

Write-Output "UseCaseCourseData Survey sheet AssessmentType column fixed."
